$wb = $excel.ActiveWorkbook

# --- 1. ResourceInfo sheet: collapse rdf_url/rdf_type/schema_url/schema_type
#        (D1:G1) into a single representation_variants column (D1). ---
$resourceInfo = $wb.Worksheets.Item("ResourceInfo")
$resourceInfo.Range("D1:G1").ClearContents()
$resourceInfo.Range("D1").Value = "representation_variants"

# --- 2. Insert a new "RepresentationVariant" sheet where "Container" used
#        to be, and push the original "Container" sheet (with its data)
#        right after it. ---
$repVariant = $wb.Worksheets.Item("Container")
$repVariant.Name = "RepresentationVariant"

$newContainer = $wb.Worksheets.Add($null, $repVariant)
$newContainer.Name = "Container"
$newContainer.Range("A1").Value = "contains_pids"

# --- 3. Populate the new RepresentationVariant sheet's header row. ---
$repVariant.Range("A1").Value = "url"
$repVariant.Range("B1").Value = "media_type"
$repVariant.Range("C1").Value = "encoding_format"
$repVariant.Range("D1").Value = "size"
